$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Update the cached TIME field result: 28-2-2018 -> 26-3-2018
# ---------------------------------------------------------------
$d.Content.Find.Execute("28-2-2018", $true, $false, $false, $false, $false, $true, 1, $false, "26-3-2018", 2) | Out-Null

# ---------------------------------------------------------------
# 2) Insert a new "Must" bullet right after the "ondersteuning" bullet
#    (before the blank paragraph / "Should" heading).
# ---------------------------------------------------------------
$anchorA = "Het design van de UI helpt hierbij, door het overzichtelijk te maken. "
$runA1 = "De gebruiker krijgt de mogelijkheid om "
$runA2 = "contact op te nemen met de beheerder via een web formulier. "
$newTextA = $anchorA + "^p" + $runA1 + "^p" + $runA2
$d.Content.Find.Execute($anchorA, $false, $false, $false, $false, $false, $true, 1, $false, $newTextA, 2) | Out-Null

# Merge the two temporary paragraphs back into one, so runA1/runA2 become
# two runs inside a single new list paragraph (instead of 2 paragraphs).
$rngA = $d.Content
$rngA.Find.Execute($runA1) | Out-Null
$rngA.Collapse(0)
$pilcrowA = $d.Range($rngA.Start, $rngA.Start + 1)
$pilcrowA.Delete()

# Make sure both new runs carry the nl-NL language tag like the rest of the doc.
$rngA1 = $d.Content
$rngA1.Find.Execute($runA1) | Out-Null
$rngA1.LanguageID = "nl-NL"

$rngA2 = $d.Content
$rngA2.Find.Execute($runA2) | Out-Null
$rngA2.LanguageID = "nl-NL"

# ---------------------------------------------------------------
# 3) Insert two new "Should" bullets right after the "Canvas" bullet.
# ---------------------------------------------------------------
$anchorB = "Het overzicht binnen Canvas wordt aangepast, zodat de gebruiker de mogelijkheid heeft om het overzicht van courses te personaliseren; Courses kunnen verwijderd worden of er kan gefilterd worden op actuele courses."
$runB1 = "De beheerder moet meerdere services kunnen toevoegen aan de bestaande portal. In dit proces kan de beheerder ook de pagina van de service vormgeven. "
$runB2 = "De beheerder kan een uitdraai vragen van alle contactmomenten tussen beheerder-gebruiker, en alle aanvragen van nieuwe services (API requests)."
$newTextB = $anchorB + "^p" + $runB1 + "^p" + $runB2
$d.Content.Find.Execute($anchorB, $false, $false, $false, $false, $false, $true, 1, $false, $newTextB, 2) | Out-Null

$rngB1 = $d.Content
$rngB1.Find.Execute($runB1) | Out-Null
$rngB1.LanguageID = "nl-NL"

$rngB2 = $d.Content
$rngB2.Find.Execute($runB2) | Out-Null
$rngB2.LanguageID = "nl-NL"
